$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3
$ws.Range("C2").Value = 4244.333333333334
$ws.Range("E2").Value = 117666.6666666667
$ws.Range("H2").Value = 5761
